$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.107.88'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.38%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.562.89'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.21%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.65'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.86'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.88%  '
$ws.Range('E7').Value = '  -2.09%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('E9').Value = '  -2.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.47'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0813'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.79'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.116'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.956.14'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.552.12'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.886'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.27'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.73%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.153.94'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.78'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0985'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.58'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.17'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.61%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '255.48'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.92%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.94'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.52%  '
$ws.Range('E25').Value = '  -5.85%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '29.10'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.07%  '
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.28'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.14'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('E30').Value = '  -5.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '153.17'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.22%  '
$ws.Range('E33').Value = '  -5.14%  '
$ws.Range('E34').Value = '  -1.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.39'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0799'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.90%  '
$ws.Range('E37').Value = '  -2.81%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '17.69'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +12.44%  '
$ws.Range('E39').Value = '  -2.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.30'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.20'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +37.77%  '
$ws.Range('E42').Value = '  -2.37%  '
$ws.Range('E43').Value = '  -1.68%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.90'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.45%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.103.49'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.96%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.18'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.09'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.811.83'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '105.80'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.92%  '
$ws.Range('E51').Value = '  +0.87%  '
